$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q5) updates
$ws.Range("B7").Value = -0.00425486747011913
$ws.Range("C7").Value = 0.3443463908930565
$ws.Range("D7").Value = 0.195059600954732
$ws.Range("E7").Value = 0.4416555229528235
$ws.Range("F7").Value = 0.4684246835828965
$ws.Range("G7").Value = 9

# Row 8 (Q6) updates
$ws.Range("B8").Value = -0.06287940768484762
$ws.Range("C8").Value = 0.2866386044041229
$ws.Range("D8").Value = 0.1396032209128983
$ws.Range("E8").Value = 0.3736351441083913
$ws.Range("F8").Value = 0.3906476367871226
$ws.Range("G8").Value = 9
